# Read data (as if from an external Excel source) and write it into the
# active workbook's active sheet, replacing the small "name/code" sample
# table with the full stock-item export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data carried an explicit "General" number format on every
# cell, so apply that to the whole target range before writing values.
$ws.Range("A1:I4").NumberFormat = "General"

$headers = @(
    "Stockcode",
    "Description",
    "Product Class",
    "Stock UOM",
    "Warehouse",
    "Supplier",
    "List Price",
    "Unit Cost",
    "Product Category"
)
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

$data = @(
    @(26012504000, "MURRPLASTIK MP 26125, R 40, LP", "MU01", "EA", 90, "INNO002", 21, 21, "M"),
    @(46012504001, "MANUALLY ADDED ITEMS",           "AB01", "EA", 90, "INNO002", 23, 21, "M"),
    @("12343-000",  "OMRON WSHR PLT OC J1 HRMNC DRV", "OM09", "EA", 30, "INNO002", 34, 12, "M")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
